$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 78, shifting existing rows 78-84 down to 79-85.
$ws.Rows.Item(78).Insert()

# Copy style of the old row 78 (now row 79) D-cell format onto new D78 so the date style persists.
$ws.Cells.Item(79, 4).Copy()
$ws.Cells.Item(78, 4).PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# Populate the new row 78 with data.
$ws.Cells.Item(78, 1).Value = 11
$ws.Cells.Item(78, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(78, 3).Value = "Bíobío"
$ws.Cells.Item(78, 4).Value = 45142
$ws.Cells.Item(78, 5).Value = 8
$ws.Cells.Item(78, 6).Value = 100112043
$ws.Cells.Item(78, 7).Value = "Pepino dulce"
$ws.Cells.Item(78, 8).Value = "Sin especificar"
$ws.Cells.Item(78, 9).Value = "Primera"
$ws.Cells.Item(78, 10).Value = 230
$ws.Cells.Item(78, 11).Value = 16000
$ws.Cells.Item(78, 12).Value = 17000
$ws.Cells.Item(78, 13).Value = 16348
$ws.Cells.Item(78, 14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(78, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(78, 16).Value = 908
$ws.Cells.Item(78, 17).Value = 18
$ws.Cells.Item(78, 18).Value = "Hortaliza"

$wb.Save()
